$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3765771687030792
$ws.Range("B1").Value = 1.112336039543152
$ws.Range("C1").Value = 6.598748683929443
$ws.Range("D1").Value = 1.839280605316162
$ws.Range("E1").Value = 1.535422801971436
